# test_set_overvote_delimiter_cvr.xlsx — switch the overvote joiner from "/" to "|"
# and pick up the accompanying whole-table white fill that shipped in the same resave.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two overvote cells ("D/A" and "A/D") become pipe-delimited ("D|A" / "A|D").
$ws.Range("C8").Value = "D|A"
$ws.Range("C9").Value = "A|D"

# The whole bordered table picked up an explicit solid white interior fill.
$ws.Range("A1:E10").Interior.Color = 16777215
